$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.452.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.625.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.52%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '324.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.05%  '

$ws.Range("E7").Value = '  -1.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.542'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.07'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.81%  '

$ws.Range("E12").Value = '  -2.49%  '

$ws.Range("E13").Value = '  +1.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.32'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.034.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.635.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.848'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.370.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.66%  '

$ws.Range("E20").Value = '  -1.50%  '

$ws.Range("E21").Value = '  -2.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0944'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '266.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.43%  '

$ws.Range("E25").Value = '  -2.58%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.74%  '

$ws.Range("E29").Value = '  -1.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.138'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.46'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0804'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.31%  '

$ws.Range("E39").Value = '  -0.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '128.33'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.82%  '

$ws.Range("E42").Value = '  -2.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0324'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.038.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.67%  '

$ws.Range("E47").Value = '  -5.36%  '

$ws.Range("E48").Value = '  -4.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.93%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.93%  '
